# Auto-generated row data and application logic for the Artfynd sheet edit.
# Rows 42-46 and 108-112 are re-ordered in place (their underlying data
# records moved position), and rows 172-174 are brand-new records appended
# at the end of the table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowData {
    param($ws, $rowNum, $data)
    foreach ($col in $data.Keys) {
        $cellInfo = $data[$col]
        $addr = "$col$rowNum"
        $type = $cellInfo[0]
        $val = $cellInfo[1]
        $cell = $ws.Range($addr)
        if ($type -eq 'b') {
            $cell.Value = [bool]([int]$val)
        } elseif ($type -eq 'n') {
            $cell.Value = [double]$val
        } else {
            if ($null -eq $val) { $val = '' }
            $cell.NumberFormat = '@'
            $cell.Value = $val
        }
    }
}

$row42 = @{
    "A" = @("n", "111896696")
    "B" = @("n", "98535")
    "C" = @("inlineStr", "Ovaliderad")
    "D" = @("inlineStr", "LC")
    "E" = @("n", "222498")
    "F" = @("inlineStr", "Blåsippa")
    "G" = @("inlineStr", "Hepatica nobilis")
    "H" = @("inlineStr", "Schreb.")
    "I" = @("inlineStr", $null)
    "P" = @("inlineStr", "Kratte masugn, Gstr")
    "Q" = @("n", "574272.5440735799")
    "R" = @("n", "6703373.755373025")
    "S" = @("n", "25")
    "T" = @("inlineStr", "Gävleborg")
    "U" = @("inlineStr", "Hofors")
    "V" = @("inlineStr", "Gästrikland")
    "W" = @("inlineStr", "Torsåker")
    "Y" = @("inlineStr", "2023-09-04")
    "Z" = @("inlineStr", "00:00")
    "AA" = @("inlineStr", "2023-09-04")
    "AB" = @("inlineStr", "00:00")
    "AD" = @("b", "0")
    "AE" = @("b", "0")
    "AG" = @("b", "0")
    "AT" = @("inlineStr", $null)
    "AW" = @("inlineStr", "Philipp Weiss")
    "AX" = @("inlineStr", "Philipp Weiss")
    "AY" = @("inlineStr", $null)
}
Set-RowData $ws 42 $row42
$ws.Range("K42").ClearContents()

$row43 = @{
    "A" = @("n", "111896703")
    "B" = @("n", "98535")
    "C" = @("inlineStr", "Ovaliderad")
    "D" = @("inlineStr", "LC")
    "E" = @("n", "222498")
    "F" = @("inlineStr", "Blåsippa")
    "G" = @("inlineStr", "Hepatica nobilis")
    "H" = @("inlineStr", "Schreb.")
    "I" = @("inlineStr", $null)
    "P" = @("inlineStr", "Kratte masugn, Gstr")
    "Q" = @("n", "574309.5943644949")
    "R" = @("n", "6703519.084753582")
    "S" = @("n", "25")
    "T" = @("inlineStr", "Gävleborg")
    "U" = @("inlineStr", "Hofors")
    "V" = @("inlineStr", "Gästrikland")
    "W" = @("inlineStr", "Torsåker")
    "Y" = @("inlineStr", "2023-09-04")
    "Z" = @("inlineStr", "00:00")
    "AA" = @("inlineStr", "2023-09-04")
    "AB" = @("inlineStr", "00:00")
    "AD" = @("b", "0")
    "AE" = @("b", "0")
    "AG" = @("b", "0")
    "AT" = @("inlineStr", $null)
    "AW" = @("inlineStr", "Philipp Weiss")
    "AX" = @("inlineStr", "Philipp Weiss")
    "AY" = @("inlineStr", $null)
}
Set-RowData $ws 43 $row43
$ws.Range("K43").ClearContents()

$row44 = @{
    "A" = @("n", "111896614")
    "B" = @("n", "90332")
    "C" = @("inlineStr", "Ovaliderad")
    "D" = @("inlineStr", "LC")
    "E" = @("n", "4769")
    "F" = @("inlineStr", "Svavelriska")
    "G" = @("inlineStr", "Lactarius scrobiculatus")
    "H" = @("inlineStr", "(Scop.:Fr.) Fr.")
    "I" = @("inlineStr", $null)
    "P" = @("inlineStr", "Kratte masugn, Gstr")
    "Q" = @("n", "574228.885558943")
    "R" = @("n", "6703430.096512586")
    "S" = @("n", "25")
    "T" = @("inlineStr", "Gävleborg")
    "U" = @("inlineStr", "Hofors")
    "V" = @("inlineStr", "Gästrikland")
    "W" = @("inlineStr", "Torsåker")
    "Y" = @("inlineStr", "2023-09-04")
    "Z" = @("inlineStr", "00:00")
    "AA" = @("inlineStr", "2023-09-04")
    "AB" = @("inlineStr", "00:00")
    "AD" = @("b", "0")
    "AE" = @("b", "0")
    "AG" = @("b", "0")
    "AT" = @("inlineStr", $null)
    "AW" = @("inlineStr", "Philipp Weiss")
    "AX" = @("inlineStr", "Philipp Weiss")
    "AY" = @("inlineStr", $null)
}
Set-RowData $ws 44 $row44

$row45 = @{
    "A" = @("n", "111886390")
    "B" = @("n", "98535")
    "C" = @("inlineStr", "Ovaliderad")
    "D" = @("inlineStr", "LC")
    "E" = @("n", "222498")
    "F" = @("inlineStr", "Blåsippa")
    "G" = @("inlineStr", "Hepatica nobilis")
    "H" = @("inlineStr", "Schreb.")
    "I" = @("inlineStr", $null)
    "K" = @("inlineStr", $null)
    "P" = @("inlineStr", "Gropbackagruvorna (Gropbackagruvorna), Gstr")
    "Q" = @("n", "574320.5008898397")
    "R" = @("n", "6703541.511550271")
    "S" = @("n", "1")
    "T" = @("inlineStr", "Gävleborg")
    "U" = @("inlineStr", "Hofors")
    "V" = @("inlineStr", "Gästrikland")
    "W" = @("inlineStr", "Torsåker")
    "Y" = @("inlineStr", "2023-09-04")
    "Z" = @("inlineStr", "11:12")
    "AA" = @("inlineStr", "2023-09-04")
    "AB" = @("inlineStr", "11:12")
    "AD" = @("b", "0")
    "AE" = @("b", "0")
    "AG" = @("b", "0")
    "AT" = @("inlineStr", $null)
    "AW" = @("inlineStr", "Annelie Hilmerby")
    "AX" = @("inlineStr", "Annelie Hilmerby")
    "AY" = @("inlineStr", $null)
}
Set-RowData $ws 45 $row45

$row46 = @{
    "A" = @("n", "111886848")
    "B" = @("n", "86223")
    "C" = @("inlineStr", "Ovaliderad")
    "D" = @("inlineStr", "NT")
    "E" = @("n", "4412")
    "F" = @("inlineStr", "Äggvaxskivling")
    "G" = @("inlineStr", "Hygrophorus karstenii")
    "H" = @("inlineStr", "Sacc. & Cub.")
    "I" = @("inlineStr", $null)
    "K" = @("inlineStr", $null)
    "P" = @("inlineStr", "Gästrikeleden (Gästrikeleden), Gstr")
    "Q" = @("n", "574236.8404598181")
    "R" = @("n", "6703596.535755535")
    "S" = @("n", "1")
    "T" = @("inlineStr", "Gävleborg")
    "U" = @("inlineStr", "Hofors")
    "V" = @("inlineStr", "Gästrikland")
    "W" = @("inlineStr", "Torsåker")
    "Y" = @("inlineStr", "2023-09-04")
    "Z" = @("inlineStr", "11:42")
    "AA" = @("inlineStr", "2023-09-04")
    "AB" = @("inlineStr", "11:42")
    "AD" = @("b", "0")
    "AE" = @("b", "0")
    "AG" = @("b", "0")
    "AT" = @("inlineStr", $null)
    "AW" = @("inlineStr", "Annelie Hilmerby")
    "AX" = @("inlineStr", "Annelie Hilmerby")
    "AY" = @("inlineStr", $null)
}
Set-RowData $ws 46 $row46

$row108 = @{
    "A" = @("n", "111885492")
    "B" = @("n", "99413")
    "C" = @("inlineStr", "Ovaliderad")
    "D" = @("inlineStr", "LC")
    "E" = @("n", "221235")
    "F" = @("inlineStr", "Vårärt")
    "G" = @("inlineStr", "Lathyrus vernus")
    "H" = @("inlineStr", "(L.) Bernh.")
    "I" = @("inlineStr", $null)
    "K" = @("inlineStr", $null)
    "P" = @("inlineStr", "Gropbackagruvorna (Gropbackagruvorna), Gstr")
    "Q" = @("n", "574346.6812743739")
    "R" = @("n", "6703445.835096321")
    "S" = @("n", "25")
    "T" = @("inlineStr", "Gävleborg")
    "U" = @("inlineStr", "Hofors")
    "V" = @("inlineStr", "Gästrikland")
    "W" = @("inlineStr", "Torsåker")
    "Y" = @("inlineStr", "2023-09-04")
    "Z" = @("inlineStr", "00:00")
    "AA" = @("inlineStr", "2023-09-04")
    "AB" = @("inlineStr", "00:00")
    "AD" = @("b", "0")
    "AE" = @("b", "0")
    "AG" = @("b", "0")
    "AT" = @("inlineStr", $null)
    "AW" = @("inlineStr", "Patric Engfeldt")
    "AX" = @("inlineStr", "Patric Engfeldt")
    "AY" = @("inlineStr", $null)
}
Set-RowData $ws 108 $row108
$ws.Range("AF108").ClearContents()

$row109 = @{
    "A" = @("n", "111885762")
    "B" = @("n", "90295")
    "C" = @("inlineStr", "Ovaliderad")
    "D" = @("inlineStr", "LC")
    "E" = @("n", "4740")
    "F" = @("inlineStr", "Sotriska")
    "G" = @("inlineStr", "Lactarius lignyotus")
    "H" = @("inlineStr", "Fr.")
    "I" = @("inlineStr", $null)
    "K" = @("inlineStr", $null)
    "P" = @("inlineStr", "Gropbackagruvorna (Gropbackagruvorna), Gstr")
    "Q" = @("n", "574311.748177869")
    "R" = @("n", "6703438.210316412")
    "S" = @("n", "25")
    "T" = @("inlineStr", "Gävleborg")
    "U" = @("inlineStr", "Hofors")
    "V" = @("inlineStr", "Gästrikland")
    "W" = @("inlineStr", "Torsåker")
    "Y" = @("inlineStr", "2023-09-04")
    "Z" = @("inlineStr", "00:00")
    "AA" = @("inlineStr", "2023-09-04")
    "AB" = @("inlineStr", "00:00")
    "AD" = @("b", "0")
    "AE" = @("b", "0")
    "AG" = @("b", "0")
    "AT" = @("inlineStr", $null)
    "AW" = @("inlineStr", "Patric Engfeldt")
    "AX" = @("inlineStr", "Patric Engfeldt")
    "AY" = @("inlineStr", $null)
}
Set-RowData $ws 109 $row109

$row110 = @{
    "A" = @("n", "111896648")
    "B" = @("n", "78512")
    "C" = @("inlineStr", "Ovaliderad")
    "D" = @("inlineStr", "LC")
    "E" = @("n", "6456")
    "F" = @("inlineStr", "Skinnlav")
    "G" = @("inlineStr", "Leptogium saturninum")
    "H" = @("inlineStr", "(Dicks.) Nyl.")
    "I" = @("inlineStr", $null)
    "P" = @("inlineStr", "Kratte masugn, Gstr")
    "Q" = @("n", "574527.773473482")
    "R" = @("n", "6703435.741965486")
    "S" = @("n", "25")
    "T" = @("inlineStr", "Gävleborg")
    "U" = @("inlineStr", "Hofors")
    "V" = @("inlineStr", "Gästrikland")
    "W" = @("inlineStr", "Torsåker")
    "Y" = @("inlineStr", "2023-09-04")
    "Z" = @("inlineStr", "00:00")
    "AA" = @("inlineStr", "2023-09-04")
    "AB" = @("inlineStr", "00:00")
    "AD" = @("b", "0")
    "AE" = @("b", "0")
    "AF" = @("inlineStr", $null)
    "AG" = @("b", "0")
    "AT" = @("inlineStr", $null)
    "AW" = @("inlineStr", "Philipp Weiss")
    "AX" = @("inlineStr", "Philipp Weiss")
    "AY" = @("inlineStr", $null)
}
Set-RowData $ws 110 $row110

$row111 = @{
    "A" = @("n", "111896717")
    "B" = @("n", "98535")
    "C" = @("inlineStr", "Ovaliderad")
    "D" = @("inlineStr", "LC")
    "E" = @("n", "222498")
    "F" = @("inlineStr", "Blåsippa")
    "G" = @("inlineStr", "Hepatica nobilis")
    "H" = @("inlineStr", "Schreb.")
    "I" = @("inlineStr", $null)
    "P" = @("inlineStr", "Kratte masugn, Gstr")
    "Q" = @("n", "574371.1760314812")
    "R" = @("n", "6703384.167737775")
    "S" = @("n", "25")
    "T" = @("inlineStr", "Gävleborg")
    "U" = @("inlineStr", "Hofors")
    "V" = @("inlineStr", "Gästrikland")
    "W" = @("inlineStr", "Torsåker")
    "Y" = @("inlineStr", "2023-09-04")
    "Z" = @("inlineStr", "00:00")
    "AA" = @("inlineStr", "2023-09-04")
    "AB" = @("inlineStr", "00:00")
    "AD" = @("b", "0")
    "AE" = @("b", "0")
    "AG" = @("b", "0")
    "AT" = @("inlineStr", $null)
    "AW" = @("inlineStr", "Philipp Weiss")
    "AX" = @("inlineStr", "Philipp Weiss")
    "AY" = @("inlineStr", $null)
}
Set-RowData $ws 111 $row111
$ws.Range("K111").ClearContents()

$row112 = @{
    "A" = @("n", "111896651")
    "B" = @("n", "89183")
    "C" = @("inlineStr", "Ovaliderad")
    "D" = @("inlineStr", "LC")
    "E" = @("n", "3215")
    "F" = @("inlineStr", "Rödgul trumpetsvamp")
    "G" = @("inlineStr", "Craterellus lutescens")
    "H" = @("inlineStr", "(Fr.) Fr.")
    "I" = @("inlineStr", $null)
    "P" = @("inlineStr", "Kratte masugn, Gstr")
    "Q" = @("n", "574442.2963542459")
    "R" = @("n", "6703458.654406736")
    "S" = @("n", "25")
    "T" = @("inlineStr", "Gävleborg")
    "U" = @("inlineStr", "Hofors")
    "V" = @("inlineStr", "Gästrikland")
    "W" = @("inlineStr", "Torsåker")
    "Y" = @("inlineStr", "2023-09-04")
    "Z" = @("inlineStr", "00:00")
    "AA" = @("inlineStr", "2023-09-04")
    "AB" = @("inlineStr", "00:00")
    "AD" = @("b", "0")
    "AE" = @("b", "0")
    "AG" = @("b", "0")
    "AT" = @("inlineStr", $null)
    "AW" = @("inlineStr", "Philipp Weiss")
    "AX" = @("inlineStr", "Philipp Weiss")
    "AY" = @("inlineStr", $null)
}
Set-RowData $ws 112 $row112
$ws.Range("K112").ClearContents()

$row172 = @{
    "A" = @("n", "112214964")
    "B" = @("n", "85313")
    "C" = @("inlineStr", "Ovaliderad")
    "D" = @("inlineStr", "NT")
    "E" = @("n", "3739")
    "F" = @("inlineStr", "Persiljespindling")
    "G" = @("inlineStr", "Cortinarius sulfurinus")
    "H" = @("inlineStr", "Quél.")
    "I" = @("inlineStr", "1")
    "J" = @("inlineStr", "fruktkroppar")
    "K" = @("inlineStr", $null)
    "N" = @("inlineStr", $null)
    "P" = @("inlineStr", "Kratte Masugn, Gstr")
    "Q" = @("n", "574526")
    "R" = @("n", "6703440")
    "S" = @("n", "25")
    "T" = @("inlineStr", "Gävleborg")
    "U" = @("inlineStr", "Hofors")
    "V" = @("inlineStr", "Gästrikland")
    "W" = @("inlineStr", "Torsåker")
    "Y" = @("inlineStr", "2023-09-19")
    "AA" = @("inlineStr", "2023-09-19")
    "AC" = @("inlineStr", "Växte tillsammans med gran och några aspar. Rikligt med dofttaggsvamp intill.")
    "AD" = @("b", "0")
    "AE" = @("b", "0")
    "AF" = @("inlineStr", $null)
    "AG" = @("b", "0")
    "AT" = @("inlineStr", $null)
    "AW" = @("inlineStr", "Henrik Tykosson")
    "AX" = @("inlineStr", "Henrik Tykosson")
    "AY" = @("inlineStr", $null)
}
Set-RowData $ws 172 $row172

$row173 = @{
    "A" = @("n", "112215073")
    "B" = @("n", "85290")
    "C" = @("inlineStr", "Ovaliderad")
    "D" = @("inlineStr", "NT")
    "E" = @("n", "6003295")
    "F" = @("inlineStr", "Odörspindling")
    "G" = @("inlineStr", "Cortinarius russeoides")
    "H" = @("inlineStr", "M.M.Moser")
    "I" = @("inlineStr", "4")
    "J" = @("inlineStr", "fruktkroppar")
    "K" = @("inlineStr", $null)
    "N" = @("inlineStr", $null)
    "P" = @("inlineStr", "Kratte Masugn, Gstr")
    "Q" = @("n", "574312")
    "R" = @("n", "6703375")
    "S" = @("n", "25")
    "T" = @("inlineStr", "Gävleborg")
    "U" = @("inlineStr", "Hofors")
    "V" = @("inlineStr", "Gästrikland")
    "W" = @("inlineStr", "Torsåker")
    "Y" = @("inlineStr", "2023-09-19")
    "AA" = @("inlineStr", "2023-09-19")
    "AC" = @("inlineStr", "Växte inom några meter ifrån Koppartaggsvamp.")
    "AD" = @("b", "0")
    "AE" = @("b", "0")
    "AF" = @("inlineStr", $null)
    "AG" = @("b", "0")
    "AT" = @("inlineStr", $null)
    "AW" = @("inlineStr", "Henrik Tykosson")
    "AX" = @("inlineStr", "Henrik Tykosson")
    "AY" = @("inlineStr", $null)
}
Set-RowData $ws 173 $row173

$row174 = @{
    "A" = @("n", "112215002")
    "B" = @("n", "87992")
    "C" = @("inlineStr", "Ovaliderad")
    "D" = @("inlineStr", "VU")
    "E" = @("n", "1593")
    "F" = @("inlineStr", "Lakritsmusseron")
    "G" = @("inlineStr", "Tricholoma apium")
    "H" = @("inlineStr", "Jul.Schäff.")
    "I" = @("inlineStr", "2")
    "J" = @("inlineStr", "fruktkroppar")
    "K" = @("inlineStr", $null)
    "N" = @("inlineStr", $null)
    "P" = @("inlineStr", "Kratte Masugn, Gstr")
    "Q" = @("n", "574226")
    "R" = @("n", "6703548")
    "S" = @("n", "100")
    "T" = @("inlineStr", "Gävleborg")
    "U" = @("inlineStr", "Hofors")
    "V" = @("inlineStr", "Gästrikland")
    "W" = @("inlineStr", "Torsåker")
    "Y" = @("inlineStr", "2023-09-19")
    "AA" = @("inlineStr", "2023-09-19")
    "AC" = @("inlineStr", "Väldigt annorlunda lokal och växtplats för en lakritsmusseron. Mestadels gran men det växte tallar i närheten också. Kalkbarrskog för övrigt och denna växte precis intill en läderdoftande fingersvamp. Fortfarande under mossan. Fruktkroppen hade den tydliga doften lakritsmusseron brukar ha.")
    "AD" = @("b", "0")
    "AE" = @("b", "0")
    "AF" = @("inlineStr", $null)
    "AG" = @("b", "0")
    "AT" = @("inlineStr", $null)
    "AW" = @("inlineStr", "Henrik Tykosson")
    "AX" = @("inlineStr", "Henrik Tykosson")
    "AY" = @("inlineStr", $null)
}
Set-RowData $ws 174 $row174
